# Greece GDP per Capita (CCode 300) — refresh the "Data" sheet's GDP-per-capita
# series (column E) with revised figures for 1833-2010 and append six new
# observations for 2011-2016 (rows 180-185), per "Update Work Week and
# Social Spending".
#
# The source values are historically stored as TEXT (shared strings), not
# numbers, so we can't just assign numeric-looking strings to .Value (Excel
# auto-converts those to the Number type). Instead we write a literal-text
# formula ( ="<text>" ) and immediately convert each formula cell to a plain
# value via Copy + PasteSpecial(xlPasteValues), which keeps the cell's type
# as Text while leaving no residual formula behind.

$ws = $excel.ActiveWorkbook.Worksheets.Item("Data")

# year|value pairs, in row order starting at row 2 (year 1833).
$data = @(
    "1833|1537",
    "1834|1635",
    "1835|2138",
    "1836|2067",
    "1837|2120",
    "1838|2230",
    "1839|2020",
    "1840|1932",
    "1841|1918",
    "1842|1959",
    "1843|1653",
    "1844|1602",
    "1845|1543",
    "1846|1691",
    "1847|1441",
    "1848|1554",
    "1849|1591",
    "1850|1607",
    "1851|1969",
    "1852|1663",
    "1853|1822",
    "1854|1908",
    "1855|2007",
    "1856|1757",
    "1857|2278",
    "1858|2197",
    "1859|2161",
    "1860|2125",
    "1861|2262",
    "1862|2283",
    "1863|2039",
    "1864|1854",
    "1865|1981",
    "1866|1890",
    "1867|1997",
    "1868|2020",
    "1869|1992",
    "1870|1938",
    "1871|1859",
    "1872|1806",
    "1873|2128",
    "1874|2013",
    "1875|1929",
    "1876|1921",
    "1877|1804",
    "1878|1985",
    "1879|2063",
    "1880|2028",
    "1881|1970",
    "1882|2098",
    "1883|2091",
    "1884|2351",
    "1885|2536",
    "1886|2426",
    "1887|2511",
    "1888|2566",
    "1889|2391",
    "1890|2257",
    "1891|1967",
    "1892|2053",
    "1893|2185",
    "1894|2233",
    "1895|2149",
    "1896|2359",
    "1897|2004",
    "1898|2165",
    "1899|2209",
    "1900|1972",
    "1901|1890",
    "1902|2423",
    "1903|2367",
    "1904|2498",
    "1905|2560",
    "1906|2530",
    "1907|2565",
    "1908|2550",
    "1909|2601",
    "1910|2592",
    "1911|3231",
    "1912|3135",
    "1913|1876",
    "1914|2394",
    "1915|1822",
    "1916|1549",
    "1917|1352",
    "1918|2279",
    "1919|2031",
    "1920|2284",
    "1921|3057",
    "1922|3129",
    "1923|3169",
    "1924|3279",
    "1925|3411",
    "1926|3475",
    "1927|3539",
    "1928|3561",
    "1929|3733",
    "1930|3599",
    "1931|3402",
    "1932|3649",
    "1933|3818",
    "1934|3854",
    "1935|3953",
    "1936|3913",
    "1937|4414",
    "1938|4267",
    "1939|4205",
    "1940|3543",
    "1941|2987",
    "1942|2517",
    "1943|2115",
    "1944|1779",
    "1945|1495",
    "1946|2209",
    "1947|2810",
    "1948|2866",
    "1949|2979",
    "1950|3052",
    "1951|3287",
    "1952|3272",
    "1953|3681",
    "1954|3759",
    "1955|4007",
    "1956|4313",
    "1957|4557",
    "1958|4723",
    "1959|4846",
    "1960|5015",
    "1961|5408",
    "1962|5577",
    "1963|6122",
    "1964|6601",
    "1965|7187",
    "1966|7570",
    "1967|7892",
    "1968|8394",
    "1969|9191",
    "1970|9900",
    "1971|10559",
    "1972|11795",
    "1973|12202",
    "1974|11716",
    "1975|12309",
    "1976|12919",
    "1977|13158",
    "1978|13860",
    "1979|14193",
    "1980|14300",
    "1981|14180",
    "1982|14153",
    "1983|14132",
    "1984|14453",
    "1985|14850",
    "1986|15047",
    "1987|14944",
    "1988|15596",
    "1989|16117",
    "1990|15964",
    "1991|16409.6894138485",
    "1992|16545.9637660561",
    "1993|16331.314704389",
    "1994|16724.2880561885",
    "1995|17148.1971986388",
    "1996|17719.5573803883",
    "1997|18588.2954879532",
    "1998|19378.7099269518",
    "1999|20076.5593750824",
    "2000|20965.3318066981",
    "2001|21913.7976809227",
    "2002|22893.9591178651",
    "2003|24380.0663273618",
    "2004|25780.3850930784",
    "2005|26091.524171943",
    "2006|27731.110923369",
    "2007|28822.915641556",
    "2008|28907.9250256727",
    "2009|27839.8975169695",
    "2010|26517.465079651",
    "2011|24349",
    "2012|22693",
    "2013|22118",
    "2014|22344",
    "2015|22442",
    "2016|22574"
)

$firstRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $firstRow + $i
    $parts = $data[$i].Split("|")
    $year = [int]$parts[0]
    $value = $parts[1]

    if ($row -gt 179) {
        # New observations (2011-2016): append full rows 180-185.
        $ws.Cells.Item($row, 1).Value = 300.0
        $ws.Cells.Item($row, 2).Value = "Greece"
        $ws.Cells.Item($row, 3).Value = "GDP per Capita"
        $ws.Cells.Item($row, 4).Value = [double]$year
    }

    # Column E ("Data") always gets the (possibly revised) text value.
    $ws.Cells.Item($row, 5).Formula = '="' + $value + '"'
}

# Flatten all the literal-text formulas down to plain text values in one
# shot (keeps t="s" shared-string cells instead of formula cells).
$lastRow = $firstRow + $data.Count - 1
$rng = $ws.Range("E" + $firstRow + ":E" + $lastRow)
$rng.Copy()
$rng.PasteSpecial(-4163)
